$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (field descriptions / help text) updates ---

# First Name / Last Name columns (D1, E1) become mandatory-for-new-user notes
$ws.Range("D1").Value = "Empty values will be ignored. For new users, this is mandatory"
$ws.Range("E1").Value = "Empty values will be ignored. For new users, this is mandatory"

# Teams column (P1): cross-reference note with Teams Managed
$ws.Range("P1").Value = "This is combined with the Teams Managed field and is not incremental, the value will be fully replaced. To clear this field, send NONE."

# Roles column (Q1): now mandatory, error instead of silent NONE-clear
$ws.Range("Q1").Value = "This field is not incremental, the value will be fully replaced. This field is mandatory - setting this to NONE will throw an error"

# Teams Managed column (R1): cross-reference note with Teams
$ws.Range("R1").Value = "This is combined with the Teams field and is not incremental, the value will be fully replaced. To clear this field, send NONE."

# --- Row 2 (field names) updates ---

# First Name / Last Name headers get a trailing "*" to flag them as mandatory
$ws.Range("D2").Value = "First Name*"
$ws.Range("E2").Value = "Last Name*"

# --- Column widths: column P/Q/R used to share one width; P and Q now need
#     to be wider to fit their longer descriptive text ---
$ws.Columns.Item(16).ColumnWidth = 135.42578125
$ws.Columns.Item(17).ColumnWidth = 111.42578125
$ws.Columns.Item(18).ColumnWidth = 80.5703125

# --- Selection: put the active cell/view on the Teams column (P2) ---
$ws.Range("P2").Select()
